# contactsImport.xlsx edit
# - Row 2 (Hannah ...): last name changes from "Weier" to "Cooper" and the
#   email cell is cleared out (including removal of its mailto hyperlink),
#   while keeping the Hyperlink cell style that was already applied to it.
# - A brand new row is inserted as row 3: another "Phil Weier" contact
#   (philweier@hotmail.com, Ninja / Vigilante / 1234567, a new note) is
#   added ahead of the pre-existing "Phil / Assasin" duplicate-contact row,
#   which is pushed down to become row 4.
# - The two mailto hyperlinks on the "Phil" rows are recreated so they keep
#   pointing at philweier@hotmail.com from their new locations (C3 and C4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlinks that currently live on C2 and C3 before we start
# moving/clearing cells around, so stale hyperlink entries don't linger.
$ws.Range("C2").Hyperlinks.Delete()
$ws.Range("C3").Hyperlinks.Delete()

# Push the existing "Phil / Assasin" row (row 3) down to row 4 by inserting
# a new blank row above it.
$ws.Rows(3).Insert()

# Populate the newly inserted row 3 with the new "Phil / Ninja" contact.
$ws.Range("A3").Value = "Phil"
$ws.Range("B3").Value = "Weier"
$ws.Range("C3").Value = "philweier@hotmail.com"
$ws.Range("D3").Value = "Ninja"
$ws.Range("E3").Value = "Vigilante"
$ws.Range("F3").Value = 1234567
$ws.Range("G3").Value = "Phil is super awesome and a huge donor."

# Update row 2: Hannah's last name becomes "Cooper" and her email is wiped.
$ws.Range("B2").Value = "Cooper"
$ws.Range("C2").Value = ""

# Re-create the mailto hyperlinks on the two Phil rows (now C4, then C3 so
# the relationship ids come out as rId1 -> C4, rId2 -> C3) and make sure
# both cells keep using the shared "Hyperlink" cell style.
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:philweier@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:philweier@hotmail.com")
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("C4").Style = "Hyperlink"

# Match the final selected cell shown in the workbook.
$ws.Range("C2").Select()
